$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","D","E","F","G","I","J","K","L")
$data = @{
    2 = @(14.59015259353513,3.516321045433621,8.894098615792716,56.0352390817772,3.789200366177493,37.94718826566889,9.318813428559046,17.70856171819532,13.52972660110104)
    3 = @(14.58269357045974,3.461793199966542,8.871404700510782,55.80099073721023,3.79288200155788,37.86491759986622,9.316441907827357,17.61328516272783,13.52871109876697)
    4 = @(14.58211000320571,3.427737810761335,8.857205343642908,55.66514502009845,3.795260186728471,37.818008908831,9.31497098677783,17.55989925620561,13.53087485395954)
    5 = @(14.58287990261657,3.413723862640665,8.851352076186492,55.61182049358995,3.796259009390832,37.79980458459255,9.314367717801884,17.53944939408213,13.53245799044744)
    6 = @(14.58306867449524,3.411388888851008,8.850376097328514,55.60308954273972,3.796426659736742,37.79683693661161,9.314267308562847,17.53613307038172,13.53276323925745)
    7 = @(14.58211630303382,3.427549353427816,8.85712667563509,55.66441759999599,3.79527353683054,37.81775970339945,9.314962866663262,17.55961815314266,13.53089336441034)
    8 = @(14.58675318981133,3.497644246145461,8.886328657572054,55.95282962171111,3.790445442822998,37.91807589592634,9.317998717181535,17.6746611195789,13.52879835640201)
    9 = @(14.62739988162841,3.630174675914309,8.941501134167314,56.58050446381853,3.781906026808002,38.14318415091343,9.323840914409852,17.93989151012775,13.54677099745335)
    10 = @(14.67624455511074,3.724074999245067,8.980779118546334,57.07770681790812,3.776191126269187,38.32555329639031,9.32807447770257,18.15751219204532,13.57337116299807)
    11 = @(14.70251716403424,3.765945703013923,8.998378625879065,57.31129513097276,3.773711157825139,38.41212394825919,9.329989660900683,18.26111892246633,13.58835973526341)
    12 = @(14.71304163904367,3.781672060206902,9.005004761926568,57.40077454782827,3.772789166161118,38.44541696476798,9.330713509178155,18.30098452689698,13.59444835231873)
    13 = @(14.71074952700301,3.778290970926686,9.003579411079588,57.38145865298323,3.772986973860198,38.4382241514209,9.330557674979186,18.29237115357794,13.59311874575496)
    14 = @(14.70337152668134,3.767242151270238,8.998924531100073,57.31863636180675,3.773634962528005,38.41485280274733,9.330049239347597,18.26438620538907,13.58885239269935)
    15 = @(14.69892702364473,3.760457387120193,8.996068284454791,57.28028808691749,3.774034100562,38.4006034074453,9.329737631192334,18.24732598036527,13.58629280187915)
    16 = @(14.6746085349,3.721320984025463,8.979623565241104,57.06258690901008,3.776355599185319,38.31996762823621,9.327949113452949,18.15083139883854,13.57244953164461)
    17 = @(14.66072309172033,3.697090083089517,8.96946667337733,56.93090337567001,3.77781036482057,38.27141892506726,9.326849282761652,18.09279423679555,13.56469514751955)
    18 = @(14.65311851000292,3.683074009068892,8.963599440379618,56.85586236964395,3.778658386902078,38.24383533550414,9.326215674583857,18.05984864953445,13.56050704005278)
    19 = @(14.65060953822931,3.678315060453547,8.961608553199449,56.83057619038879,3.778947453024505,38.23455470970707,9.326000967031424,18.04876961911041,13.55913580747645)
    20 = @(14.66216173955974,3.699677746951298,8.970550504676762,56.94484916738174,3.777654335922958,38.27655185351064,9.326966465551026,18.09892750271442,13.56549248020684)
    21 = @(14.70552306332663,3.770491022377044,9.000292823548744,57.33706131794262,3.773444168782484,38.42170374802243,9.330198615769064,18.27258916378899,13.59009434349709)
    22 = @(14.73721343530178,3.816015385568366,9.019507302247201,57.59934661479247,3.770792313351525,38.51954063775188,9.33230297947962,18.38975560782695,13.60857759220013)
    23 = @(14.71999561229963,3.791789809214015,9.009272604920751,57.45882932127929,3.772198567234009,38.46705430375704,9.331180527267358,18.32689653069998,13.5984936509189)
    24 = @(14.66151014783528,3.698508129592148,8.970060591494789,56.93854219847514,3.777724840265867,38.27423023508793,9.326913491285389,18.09615334026244,13.56513116491587)
    25 = @(14.61304910650558,3.594898206426643,8.926797596957874,56.40423128169497,3.784117493082074,38.07927817633914,9.322271409469632,17.86403457753612,13.53955032163821)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
